$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds the last-changed date for each record.
# Update every data row (2 through 34) from 2024-12-08 (serial 45634)
# to 2024-12-09 (serial 45635), keeping the existing date formatting.
for ($row = 2; $row -le 34; $row++) {
    $ws.Cells.Item($row, 3).Value = 45635
}
